# ajustes 1ra quincena mar
$wb = $excel.ActiveWorkbook

# --- mar2025: record the 1st-fortnight payments ---
$wsMar = $wb.Worksheets.Item("mar2025")

$wsMar.Range("C2").Value = 32500

$wsMar.Range("C3").Value = 65000
$wsMar.Range("C4").Value = 65000
$wsMar.Range("C5").Value = 65000
$wsMar.Range("C7").Value = 65000
$wsMar.Range("C9").Value = 65000

$wsMar.Range("C10").Value = 65000
$wsMar.Range("D10").Value = 65000

$wsMar.Range("C14").Value = 65000
$wsMar.Range("D14").Value = 65000

$wsMar.Range("C15").Value = 65000
$wsMar.Range("C17").Value = 65000
$wsMar.Range("C20").Value = 65000

$wsMar.Range("C21").Value = 50000
$wsMar.Range("D21").Value = 50000

$wsMar.Range("C24").Value = 65000

# --- update the remembered selections on both sheets ---
$wsMar.Activate()
$wsMar.Range("C10").Select()

$wsFeb = $wb.Worksheets.Item("feb2025")
$wsFeb.Activate()
$wsFeb.Range("D25").Select()
